$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - values are stored as text, so a leading
# apostrophe forces the interpreter to keep them as text rather than
# coercing to a number.
$ws.Range("D2").Value  = "'244.84"
$ws.Range("D3").Value  = "'21.92"
$ws.Range("D4").Value  = "'5.402"
$ws.Range("D5").Value  = "'0.05989"
$ws.Range("D6").Value  = "'3.392"
$ws.Range("D7").Value  = "'6.386"
$ws.Range("D8").Value  = "'0.8109"
$ws.Range("D9").Value  = "'0.9625"
$ws.Range("D11").Value = "'0.07407"
$ws.Range("D12").Value = "'0.03405"
$ws.Range("D13").Value = "'0.03061"
$ws.Range("D14").Value = "'0.09424"
$ws.Range("D16").Value = "'0.001588"
$ws.Range("D17").Value = "'0.04794"
$ws.Range("D18").Value = "'0.0005873"
$ws.Range("D19").Value = "'0.006184"
$ws.Range("D20").Value = "'0.005067"
$ws.Range("D21").Value = "'0.0009851"
$ws.Range("D23").Value = "'3.711"
$ws.Range("D26").Value = "'0.1285"
$ws.Range("D40").Value = "'0.03987"
$ws.Range("D41").Value = "'0.006593"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("D44").Value = "'0.005313"
$ws.Range("D45").Value = "'0.00005240"
$ws.Range("D47").Value = "'1.101"
$ws.Range("D48").Value = "'0.02628"

# Volume(1h) (column E) text updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
